$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 1091
$ws1.Range("F8").Value = 2137
$ws1.Range("F12").Value = 1684
$ws1.Range("F16").Value = 302
$ws1.Range("F17").Value = 224
$ws1.Range("F18").Value = 1598
$ws1.Range("F23").Value = 12353
$ws1.Range("F27").Value = 246

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 11

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 1091
$ws4.Range("F9").Value = 2137
$ws4.Range("F13").Value = 1684
$ws4.Range("F18").Value = 11
$ws4.Range("F19").Value = 302
$ws4.Range("F21").Value = 224
$ws4.Range("F27").Value = 12353
$ws4.Range("F31").Value = 246
